# Update column G ("K" - strikeouts) values for rows 2-10 on the active sheet.
# These replace the previous "Strike#" totals with actual strikeout counts.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    2  = 1
    3  = 6
    4  = 9
    5  = 5
    6  = 0
    7  = 6
    8  = 2
    9  = 1
    10 = 3
}

foreach ($row in $newValues.Keys) {
    $ws.Range("G$row").Value = $newValues[$row]
}
